# Fantasy.xlsx edit:
#   - add a new "Final" standings sheet right after "Records" (before "Schedule")
#   - move the selection on the "Records" sheet to C18
#   - populate "Final" with the end-of-season Team / Wins / Loss / PF table

$wb = $excel.ActiveWorkbook

# --- 1. Update the saved selection on the "Records" sheet -------------------
$records = $wb.Worksheets.Item("Records")
$records.Range("C18").Select()

# --- 2. Insert the new "Final" sheet, positioned before "Schedule" ----------
$schedule = $wb.Worksheets.Item("Schedule")
$final = $wb.Worksheets.Add($schedule)
$final.Name = "Final"

# --- 3. Header row ------------------------------------------------------
$final.Range("A1").Value = "Team"
$final.Range("B1").Value = "Wins"
$final.Range("C1").Value = "Loss"
$final.Range("D1").Value = "PF"

# --- 4. Final standings data ---------------------------------------------
$final.Range("A2").Value = "Assassin's Reed"
$final.Range("B2").Value = 11
$final.Range("C2").Value = 3
$final.Range("D2").Value = 1723.92

$final.Range("A3").Value = "The St. Brown Boy"
$final.Range("B3").Value = 10
$final.Range("C3").Value = 4
$final.Range("D3").Value = 1716.82

$final.Range("A4").Value = "One in Each Skibidi"
$final.Range("B4").Value = 9
$final.Range("C4").Value = 5
$final.Range("D4").Value = 1927.06

$final.Range("A5").Value = "Death to Driscoll"
$final.Range("B5").Value = 9
$final.Range("C5").Value = 5
$final.Range("D5").Value = 1619.34

$final.Range("A6").Value = "Baker Mayzyn"
$final.Range("B6").Value = 8
$final.Range("C6").Value = 6
$final.Range("D6").Value = 1886.62

$final.Range("A7").Value = "Red Wave, Red Eyes"
$final.Range("B7").Value = 8
$final.Range("C7").Value = 6
$final.Range("D7").Value = 1869.58

$final.Range("A8").Value = "Dumpster Fire"
$final.Range("B8").Value = 7
$final.Range("C8").Value = 7
$final.Range("D8").Value = 1827.66

$final.Range("A9").Value = "ElonGPT o1-mini"
$final.Range("B9").Value = 7
$final.Range("C9").Value = 7
$final.Range("D9").Value = 1722.9

# --- 5. Column A width to match the other standings sheets (~17 chars) -----
# (ColumnWidth is offset from the stored <col width> by the default ~0.833
#  padding quantum, so 16.1666... round-trips to a stored width of 17.)
$final.Columns.Item(1).ColumnWidth = 16.166666666666668

# --- 6. Leave the "Final" sheet active with D10 selected, matching the ------
#        freshly-typed-past-the-last-row cursor position.
$final.Range("D10").Select()
